# The source data set gained one additional daily observation. Insert a new
# row at position 169 (Excel shifts all rows 169-279 down to 170-280,
# extending the used range to A1:R280) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(169).Insert()

$ws.Cells.Item(169, 1).Value = 4
$ws.Cells.Item(169, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(169, 3).Value = 'Los Lagos'
$ws.Cells.Item(169, 4).Value = 44879
$ws.Cells.Item(169, 5).Value = 10
$ws.Cells.Item(169, 6).Value = 100112039
$ws.Cells.Item(169, 7).Value = 'Ciboulette'
$ws.Cells.Item(169, 8).Value = 'Sin especificar'
$ws.Cells.Item(169, 9).Value = 'Primera'
$ws.Cells.Item(169, 10).Value = 80
$ws.Cells.Item(169, 11).Value = 2500
$ws.Cells.Item(169, 12).Value = 2500
$ws.Cells.Item(169, 13).Value = 2500
$ws.Cells.Item(169, 14).Value = '$/docena de atados'
$ws.Cells.Item(169, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(169, 16).Value = 833
$ws.Cells.Item(169, 17).Value = 3
$ws.Cells.Item(169, 18).Value = 'Hortaliza'
